# cryptos list refresh (GitHub Actions scheduled data update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values that look like plain numbers need the cell
# forced to Text format first, otherwise Excel's input parser would
# silently coerce them (e.g. '563.50' -> 563.5, '0.0000180' -> 1.8E-05)
# instead of keeping the literal display string used by this sheet.

$ws.Range("D2").Value = "69.130.63"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").Value = "2.452.19"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.50"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.74"
$ws.Range("E6").Value = "  +0.71%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.509"
$ws.Range("E8").Value = "  -0.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.172"
$ws.Range("E9").Value = "  +9.12%  "

$ws.Range("E10").Value = "  -1.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.333"
$ws.Range("E11").Value = "  +1.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.63"
$ws.Range("E12").Value = "  -4.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000180"
$ws.Range("E13").Value = "  +5.69%  "

$ws.Range("D14").Value = "69.026.44"
$ws.Range("E14").Value = "  +0.73%  "

$ws.Range("D15").Value = "2.892.99"
$ws.Range("E15").Value = "  -0.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.56"
$ws.Range("E16").Value = "  +0.51%  "

$ws.Range("D17").Value = "2.445.35"
$ws.Range("E17").Value = "  -1.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.69"
$ws.Range("E18").Value = "  +1.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.56"
$ws.Range("E19").Value = "  +0.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.05"
$ws.Range("E20").Value = "  +2.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.86"
$ws.Range("E21").Value = "  +1.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.98"
$ws.Range("E22").Value = "  +5.13%  "

$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.52"
$ws.Range("E24").Value = "  -1.80%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.82"
$ws.Range("E25").Value = "  +4.53%  "

$ws.Range("D26").Value = "2.569.27"
$ws.Range("E26").Value = "  -0.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.47"
$ws.Range("E27").Value = "  +3.93%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.987"
$ws.Range("E28").Value = "  -1.21%  "

$ws.Range("D29").Value = "0.0₃0833"
$ws.Range("E29").Value = "  +1.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.23"
$ws.Range("E30").Value = "  +0.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.23"
$ws.Range("E31").Value = "  +7.30%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "437.48"
$ws.Range("E32").Value = "  +2.96%  "

$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.997"
$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.61"
$ws.Range("E34").Value = "  -0.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.28"
$ws.Range("E35").Value = "  -0.43%  "

$ws.Range("E36").Value = "  +0.16%  "

$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.108"
$ws.Range("E38").Value = "  +1.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.01"
$ws.Range("E39").Value = "  +1.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.302"
$ws.Range("E40").Value = "  +1.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.53"
$ws.Range("E41").Value = "  +3.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.41"
$ws.Range("E42").Value = "  +0.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.09"
$ws.Range("E43").Value = "  +1.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.11"
$ws.Range("E44").Value = "  +3.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.35"
$ws.Range("E45").Value = "  -0.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "130.84"
$ws.Range("E46").Value = "  +0.61%  "

$ws.Range("E47").Value = "  -0.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.486"
$ws.Range("E48").Value = "  +0.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.560"
$ws.Range("E49").Value = "  -0.21%  "

$ws.Range("E50").Value = "  +0.61%  "

$ws.Range("E51").Value = "  +2.61%  "
